$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: force a run boundary at an absolute document character position by
# toggling a no-op character formatting change on the range that starts
# there. Word (and this interop engine) splits runs at range edges when a
# formatting property is (re)applied, even if the end value is identical to
# the start value.
# ---------------------------------------------------------------------------
function Split-Boundary($pos) {
    if ($pos -le 0) { return }
    $r = $d.Range($pos, $pos + 1)
    $b = $r.Font.Bold
    $r.Font.Bold = 1
    $r.Font.Bold = $b
}

# ---------------------------------------------------------------------------
# Helper: given the start position of a paragraph run-sequence and an
# ordered list of the exact text chunks that should become separate runs,
# re-split the (already flattened/merged) paragraph text back into that
# many runs by toggling formatting at each chunk boundary.
# ---------------------------------------------------------------------------
function Apply-RunSplits($startPos, [string[]]$chunks) {
    $pos = $startPos
    for ($i = 0; $i -lt $chunks.Length; $i++) {
        if ($i -gt 0) {
            Split-Boundary $pos
        }
        $pos = $pos + $chunks[$i].Length
    }
}

$dash = [char]0x2013
$lq = [char]0x201C
$rq = [char]0x201D

# ===========================================================================
# HUNK 1: "Table names will be named ... underscores." paragraph
# ===========================================================================
$f = $d.Content
$f.Find.Execute("It has to contain one word without underscores.", $true, $false, $false, $false, $false, $true, 1, $false, "It has to contain one word without underscores signs.", 2) | Out-Null

$p = $d.Content
$p.Find.Execute("Table names will be named by class")
$pStart = $p.Start

$chunks1 = @(
    "Table names will be named by class, but in plural. It has to contain one word without underscores",
    " signs",
    "."
)
Apply-RunSplits $pStart $chunks1

# ===========================================================================
# HUNK 2: "Many-to-many related tabl|es naming pattern ..." paragraph
# ===========================================================================
$f = $d.Content
$f.Find.Execute("underscore " + $dash + " second", $true, $false, $false, $false, $false, $true, 1, $false, "underscore sign  " + $dash + " second", 2) | Out-Null

$p = $d.Content
$p.Find.Execute("Many-to-many related tabl")
$pStart = $p.Start

$chunks2 = @(
    "Many-to-many related tabl",
    "es naming pattern will be: first table name in plural " + $dash + " underscore ",
    "sign  ",
    $dash + " second table name in plural (",
    "firsttablenameinpluran_",
    " ",
    "secondtablenameinpluran",
    ")",
    "."
)
Apply-RunSplits $pStart $chunks2

# ===========================================================================
# HUNK 3: "Database columns will be ... underscores." paragraph
# ===========================================================================
$f = $d.Content
$f.Find.Execute("separated by underscores.", $true, $false, $false, $false, $false, $true, 1, $false, "separated by underscore sign s .", 2) | Out-Null

$p = $d.Content
$p.Find.Execute("Database columns will be")
$pStart = $p.Start

$chunks3 = @(
    "Database columns will be ",
    "name",
    "d in lowercase. If it is contained from more words, words will be separated by",
    " underscore",
    " sign ",
    "s",
    " ",
    "."
)
Apply-RunSplits $pStart $chunks3

# ===========================================================================
# HUNK 4: "Column with primary key ..." paragraph text change
# ===========================================================================
$f = $d.Content
$f.Find.Execute("singular " + $dash + " underscore " + $dash + " " + $lq + "id" + $rq, $true, $false, $false, $false, $false, $true, 1, $false, "singular " + $dash + " underscore sign  " + $dash + " " + $lq + "id" + $rq, 2) | Out-Null

$p = $d.Content
$p.Find.Execute("Column with primary key")
$pStart = $p.Start

$chunks4 = @(
    "Column with primary key will be named in pattern: table name in singular " + $dash + " underscore ",
    "sign  ",
    $dash + " " + $lq + "id" + $rq + ". Example: address_id. This naming rule is useful considering Hibernate mapping."
)
Apply-RunSplits $pStart $chunks4

# ===========================================================================
# HUNK 4 (continued): insert new empty paragraph + "Git" heading + body text
# right after the "Column with primary key ..." paragraph (and before the
# pre-existing trailing empty paragraphs).
# ===========================================================================
$p = $d.Content
$p.Find.Execute("This naming rule is useful considering Hibernate mapping.")
$insertPoint = $p.End

$ins = $d.Range($insertPoint, $insertPoint)
$ins.InsertParagraphAfter()

# Move to the freshly created empty paragraph and insert the "Git" heading
# paragraph after it.
$gitHeadingPos = $insertPoint + 1
$headRange = $d.Range($gitHeadingPos, $gitHeadingPos)
$headRange.InsertParagraphAfter()
$headRange.Text = "Git"
$headRange.Font.Italic = 1
$headRange.Font.Underline = 1

# Insert the body paragraph after the "Git" heading paragraph.
$bodyPos = $gitHeadingPos + 1 + 3
$bodyRange = $d.Range($bodyPos, $bodyPos)
$bodyRange.InsertParagraphAfter()
$bodyRange.Text = "Git branches will be  named in lowercase and  words in branch name will be separated by underscore sign. "
$bodyRange.Font.Italic = 0
$bodyRange.Font.Underline = 0
